# Add excel utilities (#9)
# Rename the default sheet, write a small "AddCustomerTest" sample table
# (header row + one data row) and move the selection off the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename Sheet1 -> AddCustomerTest
$ws.Name = "AddCustomerTest"

# Give the populated range an explicit (applied) cell style before writing
# values so the written cells pick up a dedicated style index.
$ws.Range("A1:C2").Style = "Normal"

# Header row
$ws.Range("A1").Value = "firstName"
$ws.Range("B1").Value = "lastName"
$ws.Range("C1").Value = "postCode"

# Data row
$ws.Range("A2").Value = "Joao"
$ws.Range("B2").Value = "Silva"
$ws.Range("C2").Value = 123456

# Move the selection to A4, below the data, matching the saved workbook's
# cursor position.
$null = $ws.Range("A4").Select()
